# Weekly update: insert a new data row for "Zanahoria" (Terminal Hortofrutícola
# Agro Chillán) at the top of its date-ordered block (row 408), pushing the
# existing rows 408-425 down to 409-426.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 408 (shifts rows 408:425 -> 409:426)
$ws.Rows.Item(408).Insert()

# Populate the newly inserted row 408 with the latest weekly price record
$ws.Cells.Item(408, 1).Value = 7
$ws.Cells.Item(408, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(408, 3).Value = "Ñuble"
$ws.Cells.Item(408, 4).Value = 45041
$ws.Cells.Item(408, 5).Value = 16
$ws.Cells.Item(408, 6).Value = 100114013
$ws.Cells.Item(408, 7).Value = "Zanahoria"
$ws.Cells.Item(408, 8).Value = "Sin especificar"
$ws.Cells.Item(408, 9).Value = "Primera"
$ws.Cells.Item(408, 10).Value = 100
$ws.Cells.Item(408, 11).Value = 7000
$ws.Cells.Item(408, 12).Value = 7000
$ws.Cells.Item(408, 13).Value = 7000
$ws.Cells.Item(408, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(408, 15).Value = "Región de Ñuble"
$ws.Cells.Item(408, 16).Value = 350
$ws.Cells.Item(408, 17).Value = 20
$ws.Cells.Item(408, 18).Value = "Hortaliza"
